$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Swap the contents of rows 83 and 84 (columns B..AB). Column A (the
#    running "id" index 81/82) stays put on each row.
# ---------------------------------------------------------------------------
$row83 = @()
$row84 = @()
for ($c = 2; $c -le 28; $c++) {
    $row83 += ,$ws.Cells.Item(83, $c).Value2
    $row84 += ,$ws.Cells.Item(84, $c).Value2
}
for ($i = 0; $i -lt $row83.Length; $i++) {
    $c = $i + 2
    $ws.Cells.Item(83, $c).Value2 = $row84[$i]
    $ws.Cells.Item(84, $c).Value2 = $row83[$i]
}

# ---------------------------------------------------------------------------
# 2) Append three new match rows (97, 98, 99) below the existing data.
#    Row 97 is a completed match (has FTHG/FTAG/FTR); rows 98 and 99 are
#    upcoming fixtures (only the opening/closing odds are known yet), whose
#    "id" values are still text placeholders rather than numeric ids.
# ---------------------------------------------------------------------------

# Copy formatting (style) from the last existing data row (96) so the new
# rows inherit the same per-column styles (bold/border on A, date format on D).
$ws.Cells.Item(96, 1).Copy() | Out-Null
$ws.Cells.Item(97, 1).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Cells.Item(98, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(99, 1).PasteSpecial(-4122) | Out-Null

$ws.Cells.Item(96, 4).Copy() | Out-Null
$ws.Cells.Item(97, 4).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(98, 4).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(99, 4).PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# --- Row 97 : York United FC 3-0 Vancouver FC (final) ---------------------
$ws.Cells.Item(97, 1).Value2 = 95
$ws.Cells.Item(97, 2).Value2 = 7802937
$ws.Cells.Item(97, 3).Value2 = "Canada Premier League"
$ws.Cells.Item(97, 4).Value2 = 45408.83333333334
$ws.Cells.Item(97, 5).Value2 = "York United FC"
$ws.Cells.Item(97, 6).Value2 = "Vancouver FC"
$ws.Cells.Item(97, 7).Value2 = 3
$ws.Cells.Item(97, 8).Value2 = 0
$ws.Cells.Item(97, 9).Value2 = "H"
$ws.Cells.Item(97, 10).Value2 = 2.2
$ws.Cells.Item(97, 11).Value2 = 3.2
$ws.Cells.Item(97, 12).Value2 = 2.9
$ws.Cells.Item(97, 13).Value2 = 2.4
$ws.Cells.Item(97, 14).Value2 = 3.3
$ws.Cells.Item(97, 15).Value2 = 2.5
$ws.Cells.Item(97, 16).Value2 = 0
$ws.Cells.Item(97, 17).Value2 = 1.85
$ws.Cells.Item(97, 18).Value2 = 1.95
$ws.Cells.Item(97, 19).Value2 = 2.75
$ws.Cells.Item(97, 20).Value2 = 2
$ws.Cells.Item(97, 21).Value2 = 1.8
$ws.Cells.Item(97, 22).Value2 = 1.4
$ws.Cells.Item(97, 23).Value2 = -1
$ws.Cells.Item(97, 24).Value2 = -1
$ws.Cells.Item(97, 25).Value2 = 0.8500000000000001
$ws.Cells.Item(97, 26).Value2 = -1
$ws.Cells.Item(97, 27).Value2 = 0.5
$ws.Cells.Item(97, 28).Value2 = -0.5

# --- Row 98 : HFX Wanderers vs Atletico Ottawa (upcoming) ------------------
$ws.Cells.Item(98, 1).Value2 = 96
$ws.Cells.Item(98, 2).NumberFormat = "@"
$ws.Cells.Item(98, 2).Value2 = "7802938"
$ws.Cells.Item(98, 2).Style = "Normal"
$ws.Cells.Item(98, 3).Value2 = "Canada Premier League"
$ws.Cells.Item(98, 4).Value2 = 45409.58333333334
$ws.Cells.Item(98, 5).Value2 = "HFX Wanderers"
$ws.Cells.Item(98, 6).Value2 = "Atletico Ottawa"
$ws.Cells.Item(98, 10).Value2 = 2
$ws.Cells.Item(98, 11).Value2 = 3.3
$ws.Cells.Item(98, 12).Value2 = 3.2
$ws.Cells.Item(98, 13).Value2 = 2.25
$ws.Cells.Item(98, 14).Value2 = 3.25
$ws.Cells.Item(98, 15).Value2 = 2.75
$ws.Cells.Item(98, 16).Value2 = -0.25
$ws.Cells.Item(98, 17).Value2 = 2.025
$ws.Cells.Item(98, 18).Value2 = 1.775
$ws.Cells.Item(98, 19).Value2 = 2.5
$ws.Cells.Item(98, 20).Value2 = 1.95
$ws.Cells.Item(98, 21).Value2 = 1.85
$ws.Cells.Item(98, 22).Value2 = 0
$ws.Cells.Item(98, 23).Value2 = 0
$ws.Cells.Item(98, 24).Value2 = 0

# --- Row 99 : Forge FC vs Valour FC (upcoming) ------------------------------
$ws.Cells.Item(99, 1).Value2 = 97
$ws.Cells.Item(99, 2).NumberFormat = "@"
$ws.Cells.Item(99, 2).Value2 = "7802876"
$ws.Cells.Item(99, 2).Style = "Normal"
$ws.Cells.Item(99, 3).Value2 = "Canada Premier League"
$ws.Cells.Item(99, 4).Value2 = 45409.70833333334
$ws.Cells.Item(99, 5).Value2 = "Forge FC"
$ws.Cells.Item(99, 6).Value2 = "Valour FC"
$ws.Cells.Item(99, 10).Value2 = 1.6
$ws.Cells.Item(99, 11).Value2 = 3.75
$ws.Cells.Item(99, 12).Value2 = 4.5
$ws.Cells.Item(99, 13).Value2 = 1.4
$ws.Cells.Item(99, 14).Value2 = 4
$ws.Cells.Item(99, 15).Value2 = 6.5
$ws.Cells.Item(99, 16).Value2 = -1.25
$ws.Cells.Item(99, 17).Value2 = 1.9
$ws.Cells.Item(99, 18).Value2 = 1.9
$ws.Cells.Item(99, 19).Value2 = 2.5
$ws.Cells.Item(99, 20).Value2 = 1.8
$ws.Cells.Item(99, 21).Value2 = 2
$ws.Cells.Item(99, 22).Value2 = 0
$ws.Cells.Item(99, 23).Value2 = 0
$ws.Cells.Item(99, 24).Value2 = 0
